$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add source/citation rows below the data table (matches the "with sources" part
# of the commit message). Write the URL cell first, then the "Source:" label,
# then the report name -- this mirrors the shared-string insertion order found
# in the target workbook.
$ws.Range("A42").Value2 = "http://www.cor.pa.gov/About%20Us/Statistics/Pages/Reports.aspx#.WU0k9evyuUk"
$ws.Range("A40").Value2 = "Source:"
$ws.Range("A40").Font.Bold = $true
$ws.Range("A41").Value2 = "Pennsylvania Department of Corrections: Annual Statistical Reports (Table 23)"

# Leave the view scrolled near the newly added rows and the cursor parked just
# past the data, matching where the editor ended up after adding the source
# information.
$ws.Range("P37").Select() | Out-Null

# Page orientation was set to portrait explicitly.
$ws.PageSetup.Orientation = 1
